$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header) updates ---
$ws.Range("AB1").Value = "NoOfManualTrigger"
$ws.Range("AC1").Value = "ExpectedRecordLength"

# --- Row 2 (data) updates ---
# Leading apostrophe keeps these as text (quote-prefixed), matching the
# existing "number-as-text" style already used in this row.
$ws.Range("L2").Value = "'2000"
$ws.Range("O2").Value = "'15"
$ws.Range("P2").Value = "'1"
$ws.Range("R2").Value = "'5000"
$ws.Range("U2").Value = "'15"
$ws.Range("V2").Value = "'1"

# New trailing columns on row 2
$ws.Range("AB2").Value = "'30"
$ws.Range("AC2").Value = "'30000"

# Auto-fit the newly populated column so its width reflects the new content
$ws.Columns.Item(28).AutoFit()

# --- Remove the now-duplicate third row entirely ---
$ws.Rows(3).Delete()

# Reproduce the author's final selection (selecting the deleted row's
# former position, i.e. the new row 3 = first empty row) before saving.
$ws.Rows(3).Select()
